$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "68.370.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -2.97%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.705.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -3.52%  "

# Row 4
$ws.Range("E4").Value2 = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "603.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +1.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "181.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +8.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "3.698.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -3.63%  "

# Row 8
$ws.Range("E8").Value2 = "  -5.89%  "

# Row 9
$ws.Range("E9").Value2 = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -4.00%  "

# Row 11
$ws.Range("E11").Value2 = "  -6.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "56.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +6.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -8.74%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "10.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -8.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "4.295.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -3.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "3.709.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -3.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "19.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -8.07%  "

# Row 18
$ws.Range("E18").Value2 = "  -2.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -6.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "1.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -6.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "68.207.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -3.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "409.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -6.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "4.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -1.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "89.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -4.89%  "

# Row 25
$ws.Range("E25").Value2 = "  -7.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "12.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -7.11%  "

# Row 27
$ws.Range("B27").Value2 = "Toncoin"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "3.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -3.33%  "

# Row 28
$ws.Range("B28").Value2 = "RenderToken"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "10.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -3.45%  "

# Row 29
$ws.Range("E29").Value2 = "  +1.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "9.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -8.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "32.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -6.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "7.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -10.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "12.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -6.97%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -6.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "43.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -8.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "64.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -7.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "600.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -4.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.0₃0890"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -8.58%  "

# Row 39
$ws.Range("B39").Value2 = "Dai"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +0.25%  "

# Row 40
$ws.Range("B40").Value2 = "TheGraph"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.399"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -5.37%  "

# Row 41
$ws.Range("E41").Value2 = "  +0.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.136"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -6.06%  "

# Row 43
$ws.Range("B43").Value2 = "ThetaToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "3.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -6.96%  "

# Row 44
$ws.Range("B44").Value2 = "Fetch.AI"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "2.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +1.62%  "

# Row 45
$ws.Range("E45").Value2 = "  -5.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -13.13%  "

# Row 47
$ws.Range("E47").Value2 = "  -7.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -3.58%  "

# Row 49
$ws.Range("B49").Value2 = "Maker"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "2.780.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -1.68%  "

# Row 50
$ws.Range("B50").Value2 = "Stellar"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.134"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -6.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "3.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -2.30%  "

